# Update the indirect method output for broker alias "hait_ehfz": append a
# new block of rows for DataDate 20210128, mirroring the BrokerAlias
# ordering (hait, huat, swhy, gtja, zx) used in the prior date block, with
# the freshly-pulled IntersectionSecurityCount results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($cell, $value) {
    # Force the cell to be written as text (shared string), matching how
    # the existing DataDate / BrokerAlias / sentinel columns are stored,
    # while avoiding any lingering explicit cell style on the result.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$newRows = @(
    @{ Row = 37; Date = "20210128"; Broker = "hait"; Count = "DataFileNotExists" },
    @{ Row = 38; Date = "20210128"; Broker = "huat"; Count = 0 },
    @{ Row = 39; Date = "20210128"; Broker = "swhy"; Count = "DataFileNotExists" },
    @{ Row = 40; Date = "20210128"; Broker = "gtja"; Count = "DataFileNotExists" },
    @{ Row = 41; Date = "20210128"; Broker = "zx";   Count = "DataFileNotExists" }
)

foreach ($r in $newRows) {
    Set-TextCell "A$($r.Row)" $r.Date
    Set-TextCell "B$($r.Row)" $r.Broker

    if ($r.Count -is [string]) {
        Set-TextCell "C$($r.Row)" $r.Count
    } else {
        $ws.Range("C$($r.Row)").Value = $r.Count
    }
}
